$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default (unstyled) cell format from column C, which is never edited,
# so it reliably reflects the workbook's base style for data rows.
$defaultStyle = $ws.Cells.Item(2, 3).Style

# Rows whose Price text looks like a plain number and would otherwise be
# auto-converted to a numeric value by Excel; force them to stay text.
$textRows = @(4, 5, 6, 10, 11, 12, 13, 14, 19, 22, 23, 24, 25, 26, 28, 29, 31, 32, 33, 36, 38, 39, 40, 41, 44, 47, 48, 49, 50, 51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '67.695.35'
$ws.Cells.Item(2, 5).Value = '  +5.77%  '
$ws.Cells.Item(3, 4).Value = '3.495.74'
$ws.Cells.Item(3, 5).Value = '  +6.20%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.24%  '
$ws.Cells.Item(5, 4).Value = '189.60'
$ws.Cells.Item(5, 5).Value = '  +9.44%  '
$ws.Cells.Item(6, 4).Value = '555.45'
$ws.Cells.Item(6, 5).Value = '  +6.31%  '
$ws.Cells.Item(7, 5).Value = '  +2.02%  '
$ws.Cells.Item(8, 4).Value = '3.486.54'
$ws.Cells.Item(8, 5).Value = '  +6.18%  '
$ws.Cells.Item(9, 5).Value = '  -0.13%  '
$ws.Cells.Item(10, 4).Value = '0.640'
$ws.Cells.Item(10, 5).Value = '  +5.92%  '
$ws.Cells.Item(11, 4).Value = '57.01'
$ws.Cells.Item(11, 5).Value = '  +1.64%  '
$ws.Cells.Item(12, 4).Value = '0.150'
$ws.Cells.Item(12, 5).Value = '  +12.99%  '
$ws.Cells.Item(13, 4).Value = '0.0000276'
$ws.Cells.Item(14, 4).Value = '9.49'
$ws.Cells.Item(14, 5).Value = '  +5.15%  '
$ws.Cells.Item(15, 4).Value = '4.045.88'
$ws.Cells.Item(15, 5).Value = '  +5.95%  '
$ws.Cells.Item(16, 4).Value = '3.493.90'
$ws.Cells.Item(16, 5).Value = '  +6.20%  '
$ws.Cells.Item(17, 4).Value = '68.020.08'
$ws.Cells.Item(17, 5).Value = '  +6.36%  '
$ws.Cells.Item(18, 5).Value = '  +4.39%  '
$ws.Cells.Item(19, 4).Value = '18.39'
$ws.Cells.Item(19, 5).Value = '  +6.06%  '
$ws.Cells.Item(20, 5).Value = '  +7.62%  '
$ws.Cells.Item(21, 5).Value = '  +6.01%  '
$ws.Cells.Item(22, 4).Value = '404.88'
$ws.Cells.Item(22, 5).Value = '  +8.57%  '
$ws.Cells.Item(23, 4).Value = '12.18'
$ws.Cells.Item(23, 5).Value = '  +11.10%  '
$ws.Cells.Item(24, 4).Value = '3.97'
$ws.Cells.Item(24, 5).Value = '  +5.77%  '
$ws.Cells.Item(25, 4).Value = '84.79'
$ws.Cells.Item(25, 5).Value = '  +5.97%  '
$ws.Cells.Item(26, 4).Value = '4.22'
$ws.Cells.Item(26, 5).Value = '  +7.72%  '
$ws.Cells.Item(27, 5).Value = '  +9.00%  '
$ws.Cells.Item(28, 4).Value = '6.25'
$ws.Cells.Item(28, 5).Value = '  +2.81%  '
$ws.Cells.Item(29, 4).Value = '11.88'
$ws.Cells.Item(29, 5).Value = '  +4.95%  '
$ws.Cells.Item(30, 5).Value = '  +4.60%  '
$ws.Cells.Item(31, 4).Value = '30.40'
$ws.Cells.Item(31, 5).Value = '  +6.00%  '
$ws.Cells.Item(32, 4).Value = '685.73'
$ws.Cells.Item(32, 5).Value = '  +7.15%  '
$ws.Cells.Item(33, 4).Value = '6.93'
$ws.Cells.Item(33, 5).Value = '  +5.19%  '
$ws.Cells.Item(34, 5).Value = '  +4.60%  '
$ws.Cells.Item(35, 5).Value = '  +5.74%  '
$ws.Cells.Item(36, 4).Value = '59.94'
$ws.Cells.Item(36, 5).Value = '  +2.00%  '
$ws.Cells.Item(37, 4).Value = '0.0₃0842'
$ws.Cells.Item(37, 5).Value = '  +21.39%  '
$ws.Cells.Item(38, 4).Value = '39.11'
$ws.Cells.Item(38, 5).Value = '  +6.96%  '
$ws.Cells.Item(39, 4).Value = '0.405'
$ws.Cells.Item(39, 5).Value = '  +5.84%  '
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  +0.00%  '
$ws.Cells.Item(41, 4).Value = '3.44'
$ws.Cells.Item(41, 5).Value = '  +25.26%  '
$ws.Cells.Item(42, 5).Value = '  +15.42%  '
$ws.Cells.Item(43, 5).Value = '  +10.66%  '
$ws.Cells.Item(44, 4).Value = '0.999'
$ws.Cells.Item(44, 5).Value = '  +0.06%  '
$ws.Cells.Item(45, 4).Value = '3.047.30'
$ws.Cells.Item(45, 5).Value = '  +4.09%  '
$ws.Cells.Item(46, 5).Value = '  +10.74%  '
$ws.Cells.Item(47, 4).Value = '0.0423'
$ws.Cells.Item(47, 5).Value = '  +6.82%  '
$ws.Cells.Item(48, 4).Value = '3.27'
$ws.Cells.Item(48, 5).Value = '  +10.70%  '
$ws.Cells.Item(49, 4).Value = '2.75'
$ws.Cells.Item(49, 5).Value = '  +4.41%  '
$ws.Cells.Item(50, 4).Value = '8.99'
$ws.Cells.Item(50, 5).Value = '  +15.66%  '
$ws.Cells.Item(51, 4).Value = '0.130'
$ws.Cells.Item(51, 5).Value = '  +4.51%  '

# Restore the original (General) style on the forced-text cells so only the
# cell values change, matching the source formatting.
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).Style = $defaultStyle
}
